# Revert "Adding the RES Hourly Production Forecast to the Portfolio"
#
# 1) The "Lookup" column (D) text cached in the shared-string table is built
#    from a date prefix ("dd.mm.yyyy") + the Interval number (column B).
#    The date prefix moves from 24.09.2024 back to 29.08.2024.
# 2) The "Interval" timestamps in column A move back 26 days
#    (45559.x -> 45533.x, i.e. 24.09.2024 -> 29.08.2024).
# 3) The "Prediction" values in column C (rows 31-82) revert to the older
#    forecast numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the cached "dd.mm.yyyy" prefix used to build the Lookup text.
# Using Replace keeps the shared-string indices stable (in-place text swap).
$ws.Cells.Replace("24.09.2024", "29.08.2024")

# --- 2) Shift every Interval timestamp in column A back by 26 days.
$lastRow = 96
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2()
    $cell.Value = $serial - 26
}

# --- 3) Restore the previous Prediction values for rows 31-82.
$predictions = @{
    31 = 0.011
    32 = 0.015
    33 = 0.028
    34 = 0.07099999999999999
    35 = 0.137
    36 = 0.204
    37 = 0.26
    38 = 0.297
    39 = 0.354
    40 = 0.397
    41 = 0.378
    42 = 0.45
    43 = 0.492
    44 = 0.515
    45 = 0.529
    46 = 0.5649999999999999
    47 = 0.594
    48 = 0.613
    49 = 0.627
    50 = 0.629
    51 = 0.629
    52 = 0.629
    53 = 0.623
    54 = 0.617
    55 = 0.598
    56 = 0.592
    57 = 0.587
    58 = 0.583
    59 = 0.5620000000000001
    60 = 0.514
    61 = 0.504
    62 = 0.476
    63 = 0.455
    64 = 0.45
    65 = 0.446
    66 = 0.426
    67 = 0.403
    68 = 0.337
    69 = 0.282
    70 = 0.247
    71 = 0.213
    72 = 0.174
    73 = 0.143
    74 = 0.118
    75 = 0.08799999999999999
    76 = 0.068
    77 = 0.058
    78 = 0.045
    79 = 0.039
    80 = 0.03
    81 = 0.021
    82 = 0.013
}

foreach ($row in $predictions.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $predictions[$row]
}
